$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that swap values between row 2 and row 5: D, L, M, N, O, P, S
$cols = @("D", "L", "M", "N", "O", "P", "S")

foreach ($col in $cols) {
    $addr2 = $col + "2"
    $addr5 = $col + "5"
    $val2 = $ws.Range($addr2).Value2
    $val5 = $ws.Range($addr5).Value2
    $ws.Range($addr2).Value = $val5
    $ws.Range($addr5).Value = $val2
}
